$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.203.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.863.80"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06536"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07850"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.37"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.870.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.099"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.202.19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.522"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.64"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.116.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007284"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.191"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.15"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.68%  "
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.920"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.374"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09698"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.412"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.474"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.087"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04686"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.115"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7042"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.728"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01847"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.257"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.78%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.533"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.50"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8458"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.06"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.171"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.170"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "935.18"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.67%  "
$ws.Range("E50").Value = "  -0.42%  "
